# Rename the "LOT" (pop_filter2) and "Subpopulation" (pop_filter1) page-name
# cells on Sheet1, and update the selected/visible cell to match the new
# author view, per commit "Updated the page names for Subpopulation and LOT
# pages".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 block: lot_section* -> pop_filter2_section*
# (write in I, G, H order so new shared-string entries land in the same
# order they do in the target workbook)
$ws.Range("I2").Value2 = "pop_filter2_section"
$ws.Range("G2").Value2 = "pop_filter2_section1"
$ws.Range("H2").Value2 = "pop_filter2_section1_checkbox"

# Row 7 block: sub_pop_section* -> pop_filter1_section*
$ws.Range("I7").Value2 = "pop_filter1_section"
$ws.Range("H7").Value2 = "pop_filter1_section1_checkbox"
$ws.Range("G7").Value2 = "pop_filter1_section1"

# Update the view: scroll to D1 (topLeftCell) and select I11.
$ws.Range("I11").Select()
